$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New bug row (row 5) -------------------------------------------------
$ws.Range("A5").Value = "Cursor not rotating"
$ws.Range("B5").Value = "PR-HUD"
$ws.Range("C5").Value = "Genio"

# --- New columns: SOLVED? (E) and REMARKS (F) ----------------------------
$ws.Range("E1").Value = "SOLVED?"
$ws.Range("E2").Value = "NO"
$ws.Range("E3").Value = "NO"
$ws.Range("E4").Value = "NO"
$ws.Range("E5").Value = "NO"
$ws.Range("F1").Value = "REMARKS"

# --- Borders: reshuffle header border styles now that F is last column --
# F1 takes over the "last column" (right edge) border that D1 used to have
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
# D1 and E1 become interior header cells, matching B1/C1's border
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# --- Apply the built-in "Bad" cell style to the SOLVED?=NO column --------
$ws.Range("E2:E5").Style = "Bad"

# --- Column widths ---------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 24.053385416666668
$ws.Columns.Item(3).ColumnWidth = 14.608072916666666
$ws.Columns.Item(6).ColumnWidth = 40.721354166666664

# --- Selection shown when the file is reopened -----------------------------
$ws.Range("B18").Select()
